$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "25+58=",
    "63+7=",
    "43+28=",
    "62-27=",
    "13+64=",
    "43+49=",
    "49-6=",
    "85-5=",
    "48-9=",
    "37-25=",
    "80-73=",
    "62-30=",
    "47-29=",
    "65+26=",
    "11+34=",
    "5+47=",
    "45+18=",
    "42-21=",
    "44-4=",
    "86-81=",
    "0+0=",
    "96-59=",
    "81-27=",
    "73+11=",
    "37+10=",
    "39+3=",
    "21+59=",
    "38+50=",
    "65-11=",
    "5+91=",
    "29+40=",
    "42-35=",
    "58-4=",
    "24+27=",
    "53-19=",
    "81-46=",
    "87-29=",
    "61+1=",
    "37+21=",
    "58-11=",
    "54+23=",
    "51-20=",
    "14+12=",
    "35+62=",
    "79+4=",
    "39-11=",
    "42-1=",
    "99-43=",
    "3+9=",
    "78+13=",
    "28+16=",
    "83-57=",
    "90-9=",
    "78-20=",
    "40+50=",
    "19+29=",
    "66-27=",
    "35+11=",
    "91-25=",
    "91-89=",
    "85+7=",
    "22+47=",
    "81+11=",
    "67-40=",
    "12+53=",
    "91-50=",
    "98-37=",
    "28+68=",
    "74-64=",
    "56+41=",
    "41+42=",
    "38+14=",
    "22-1=",
    "17+16=",
    "30+14=",
    "85-57=",
    "90-1=",
    "57-53=",
    "79-15=",
    "47-7=",
    "43+44=",
    "70-67=",
    "4+14=",
    "79-2=",
    "94-9=",
    "32+20=",
    "53-28=",
    "44+37=",
    "35-23=",
    "4+12=",
    "90-81=",
    "48-34=",
    "35+63=",
    "72-15=",
    "9+84=",
    "61-47=",
    "54+30=",
    "95-46=",
    "34+48=",
    "39-27="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Done, wrote" $idx "cells"
